$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (a newer quarter + an additional
# quarter were added to each of the three tables on the sheet), which
# shifts the old D:K data right to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# The newly inserted columns pick up column C's formatting by default;
# copy the number formats from column F (the old column D, now shifted)
# onto the two new columns so dates/numbers render the same as their
# neighbours.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new leading columns (and the handful of existing
# columns on row 91 whose historical figures were corrected) with the
# updated financial figures.
  $ws.Range("D7").Value = 43465; $ws.Range("E7").Value = 43373
  $ws.Range("D8").Value = 137900; $ws.Range("E8").Value = 101200
  $ws.Range("D9").Value = 3600; $ws.Range("E9").Value = 3800
  $ws.Range("D10").Value = 134300; $ws.Range("E10").Value = 97400
  $ws.Range("D12").Value = "NA"; $ws.Range("E12").Value = "NA"
  $ws.Range("D13").Value = 0; $ws.Range("E13").Value = 0
  $ws.Range("D14").Value = 2200; $ws.Range("E14").Value = 1500
  $ws.Range("D15").Value = 39500; $ws.Range("E15").Value = 34400
  $ws.Range("D17").Value = 139000; $ws.Range("E17").Value = 100500
  $ws.Range("D18").Value = -1100; $ws.Range("E18").Value = 700
  $ws.Range("D20").Value = 200; $ws.Range("E20").Value = 600
  $ws.Range("D21").Value = 47700; $ws.Range("E21").Value = 40400
  $ws.Range("D22").Value = 0; $ws.Range("E22").Value = 0
  $ws.Range("D23").Value = -900; $ws.Range("E23").Value = 1300
  $ws.Range("D24").Value = -200; $ws.Range("E24").Value = 600
  $ws.Range("D25").Value = 0; $ws.Range("E25").Value = 0
  $ws.Range("D26").Value = -600; $ws.Range("E26").Value = 700
  $ws.Range("D27").Value = 1000; $ws.Range("E27").Value = 4600
  $ws.Range("D28").Value = 0; $ws.Range("E28").Value = 0
  $ws.Range("D29").Value = "NA"; $ws.Range("E29").Value = "NA"
  $ws.Range("D30").Value = 0; $ws.Range("E30").Value = 0
  $ws.Range("D31").Value = 0; $ws.Range("E31").Value = 0
  $ws.Range("D32").Value = -200; $ws.Range("E32").Value = -600
  $ws.Range("D33").Value = 1000; $ws.Range("E33").Value = 4600
  $ws.Range("D34").Value = 0; $ws.Range("E34").Value = 0
  $ws.Range("D35").Value = 1000; $ws.Range("E35").Value = 4600
  $ws.Range("D38").Value = 43465; $ws.Range("E38").Value = 43373
  $ws.Range("D41").Value = 99600; $ws.Range("E41").Value = 163100
  $ws.Range("D42").Value = 0; $ws.Range("E42").Value = 0
  $ws.Range("D43").Value = 53800; $ws.Range("E43").Value = 50000
  $ws.Range("D44").Value = 0; $ws.Range("E44").Value = 0
  $ws.Range("D45").Value = 21200; $ws.Range("E45").Value = 22400
  $ws.Range("D46").Value = 174600; $ws.Range("E46").Value = 235500
  $ws.Range("D47").Value = 59200; $ws.Range("E47").Value = 59100
  $ws.Range("D48").Value = 2141100; $ws.Range("E48").Value = 1943500
  $ws.Range("D49").Value = 155100; $ws.Range("E49").Value = 152200
  $ws.Range("D50").Value = 0; $ws.Range("E50").Value = 0
  $ws.Range("D51").Value = 0; $ws.Range("E51").Value = 0
  $ws.Range("D52").Value = 108800; $ws.Range("E52").Value = 82800
  $ws.Range("D53").Value = 0; $ws.Range("E53").Value = 0
  $ws.Range("D54").Value = 2638800; $ws.Range("E54").Value = 2473200
  $ws.Range("D57").Value = 112200; $ws.Range("E57").Value = 83600
  $ws.Range("D58").Value = 71700; $ws.Range("E58").Value = 72900
  $ws.Range("D59").Value = 38500; $ws.Range("E59").Value = 33900
  $ws.Range("D60").Value = 222400; $ws.Range("E60").Value = 190500
  $ws.Range("D61").Value = 1165700; $ws.Range("E61").Value = 1065000
  $ws.Range("D62").Value = 196900; $ws.Range("E62").Value = 155400
  $ws.Range("D63").Value = 0; $ws.Range("E63").Value = 0
  $ws.Range("D64").Value = 0; $ws.Range("E64").Value = 0
  $ws.Range("D65").Value = 0; $ws.Range("E65").Value = 0
  $ws.Range("D66").Value = 1641400; $ws.Range("E66").Value = 1468700
  $ws.Range("D68").Value = 0; $ws.Range("E68").Value = 0
  $ws.Range("D69").Value = 0; $ws.Range("E69").Value = 0
  $ws.Range("D70").Value = 0; $ws.Range("E70").Value = 0
  $ws.Range("D71").Value = 0; $ws.Range("E71").Value = 0
  $ws.Range("D72").Value = -32800; $ws.Range("E72").Value = -33900
  $ws.Range("D73").Value = 0; $ws.Range("E73").Value = 0
  $ws.Range("D74").Value = 0; $ws.Range("E74").Value = 0
  $ws.Range("D75").Value = 0; $ws.Range("E75").Value = 0
  $ws.Range("D76").Value = 997400; $ws.Range("E76").Value = 1004500
  $ws.Range("D77").Value = 0; $ws.Range("E77").Value = 0
  $ws.Range("D80").Value = 43465; $ws.Range("E80").Value = 43373
  $ws.Range("D81").Value = 1000; $ws.Range("E81").Value = 4600
  $ws.Range("D83").Value = 48500; $ws.Range("E83").Value = 39100
  $ws.Range("D84").Value = 0; $ws.Range("E84").Value = 0
  $ws.Range("D85").Value = 0; $ws.Range("E85").Value = 0
  $ws.Range("D86").Value = 0; $ws.Range("E86").Value = 0
  $ws.Range("D87").Value = 0; $ws.Range("E87").Value = 0
  $ws.Range("D88").Value = 0; $ws.Range("E88").Value = 0
  $ws.Range("D89").Value = 47300; $ws.Range("E89").Value = 27300
  $ws.Range("D91").Value = 7200; $ws.Range("E91").Value = 500
  $ws.Range("D92").Value = 0; $ws.Range("E92").Value = 0
  $ws.Range("D93").Value = 0; $ws.Range("E93").Value = 0
  $ws.Range("D94").Value = -201900; $ws.Range("E94").Value = -176500
  $ws.Range("D96").Value = -28000; $ws.Range("E96").Value = -28000
  $ws.Range("D97").Value = 0; $ws.Range("E97").Value = 0
  $ws.Range("D98").Value = 0; $ws.Range("E98").Value = 0
  $ws.Range("D99").Value = 0; $ws.Range("E99").Value = 0
  $ws.Range("D100").Value = 90000; $ws.Range("E100").Value = 260900
  $ws.Range("D101").Value = 0; $ws.Range("E101").Value = 0
  $ws.Range("D102").Value = -64700; $ws.Range("E102").Value = 111700

# Row 91 ("Changes In Other Operating Activities") also had its D:J
# figures corrected as part of this update, not merely shifted.
$ws.Range("D91").Value = 7200
$ws.Range("E91").Value = 500
$ws.Range("F91").Value = -11000
$ws.Range("G91").Value = -6900
$ws.Range("H91").Value = -500
$ws.Range("I91").Value = -6100
$ws.Range("J91").Value = -4600
